$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 'ACORN Participating Countries'
$ws.Range("A12").Value = 'All ''orgname'' are provided.'
$ws.Range("A13").Value = 'All ''patid'' are provided.'
$ws.Range("A14").Value = 'All ''specdate'' are provided.'
$ws.Range("A15").Value = 'All ''specdate'' are today or before today.'
$ws.Range("A16").Value = 'All ''specgroup'' are provided.'
$ws.Range("A17").Value = 'All ''specid'' are provided.'
$ws.Range("A18").Value = 'All dates of enrolment for HAI patients have a matching date in the HAI survey dataset'
$ws.Range("A19").Value = 'All Other Organisms'
$ws.Range("A20").Value = 'All valid records have an ACORN ID.'
$ws.Range("A21").Value = 'AMR'
$ws.Range("A22").Value = 'and generate enrolment log.'
$ws.Range("A23").Value = 'Attempting to connect.'
$ws.Range("A24").Value = 'Blood culture collected within 24 hours of admission (CAI) / symptom onset (HAI)'
$ws.Range("A25").Value = 'Blood Culture Contaminants'
$ws.Range("A26").Value = 'Bloodstream Infection (BSI)'
$ws.Range("A27").Value = 'Calculated age is consistent with ''Age Category'''
$ws.Range("A28").Value = 'Calculated age isn''t always consistent with ''Age Category'''
$ws.Range("A29").Value = 'Cancel'
$ws.Range("A30").Value = 'Care should be taken when interpreting rates and AMR profiles where there are small numbers of cases or bacterial isolates: point estimates may be unreliable.'
$ws.Range("A31").Value = 'Clinical and day-28 outcomes are consistent.'
$ws.Range("A32").Value = 'Clinical and day-28 outcomes aren''t consistent for some dead patients.'
$ws.Range("A33").Value = 'Clinical Outcome'
$ws.Range("A34").Value = 'Clinical Outcome Status:'
$ws.Range("A35").Value = 'Co-resistances'
$ws.Range("A36").Value = 'Combine Susceptible + Intermediate'
$ws.Range("A37").Value = 'Consider saving .acorn file on the cloud for additional security.'
$ws.Range("A38").Value = 'Contains names of organisms before and after mapping.'
$ws.Range("A39").Value = 'Couldn''t connect to server. Please check internet access.'
$ws.Range("A40").Value = 'Critical errors with clinical data.'
$ws.Range("A41").Value = 'Culture results per specimen type'
$ws.Range("A42").Value = 'Data Management'
$ws.Range("A43").Value = 'Date of Enrolment'
$ws.Range("A44").Value = 'Day 28'
$ws.Range("A45").Value = 'Day 28 Status:'
$ws.Range("A46").Value = 'Diagnosis at Enrolment'
$ws.Range("A47").Value = 'Dismiss'
$ws.Range("A48").Value = 'Distribution of Enrolments'
$ws.Range("A49").Value = 'Download Enrolment Log (.xlsx)'
$ws.Range("A50").Value = 'Download Lab Log (.xlsx)'
$ws.Range("A70").Value = 'HAI point prevalence by '
$ws.Range("A110").Value = 'Remove ''Not Cultured'' specimens'
$ws.Range("A111").Value = 'Remove blood culture contaminants from the following visualizations'
$ws.Range("A112").Value = 'Reset Enrolments Filters'
$ws.Range("A113").Value = 'Resistance to 3rd gen. Cephalosporins Over Time'
$ws.Range("A114").Value = 'Resistance to Carbapenems Over Time'
$ws.Range("A115").Value = 'Resistance to Fluoroquinolones Over Time'
$ws.Range("A116").Value = 'Resistance to Oxacillin Over Time'
$ws.Range("A117").Value = 'Resistance to Penicillin G - meningitis Over Time'
$ws.Range("A118").Value = 'Resistance to Penicillin G Over Time'
$ws.Range("A119").Value = 'Retriving data from REDCap server.'
$ws.Range("A120").Value = 'Save .acorn file'
$ws.Range("A121").Value = 'Save acorn data'
$ws.Range("A122").Value = 'Save on Server'
$ws.Range("A123").Value = 'See Breakdown by Ward'
$ws.Range("A124").Value = 'See by Week'
